$wb = $excel.ActiveWorkbook

$calcs = $wb.Worksheets.Item("Calcs")
$calcs.Activate()
$calcs.Range("B123").Select()

$inputs = $wb.Worksheets.Item("Inputs")
$inputs.Activate()
